$d = $word.ActiveDocument

# 1. <move> element -> <event> element
$d.Content.Find.Execute("move", $true, $false, $false, $false, $false, $true, 1, $false, "event", 2)

# 2. @who for action-doer -> @resp for action-doer
$d.Content.Find.Execute("who for action-doer", $true, $false, $false, $false, $false, $true, 1, $false, "resp for action-doer", 2)

# 3. @type -> @style (the one in the actions section, after the blank paragraph following "for action-doer")
$d.Content.Find.Execute("type", $true, $false, $false, $false, $false, $true, 1, $false, "style", 2)
